$d = $word.ActiveDocument

# --- 1. Remove the 5 paragraphs describing the "Зачислен на учебу ..." block,
#        through the "Проживает в общежитии ... {numRental}" paragraph.
#        These are located right after the paragraph containing "{numRoom}"
#        and right before "Согласовано в установленном порядке."
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "Зачислен на учебу") {
        $startPara = $i
    }
    if ($t -match "numRental") {
        $endPara = $i
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $delRange = $d.Range($d.Paragraphs.Item($startPara).Range.Start, $d.Paragraphs.Item($endPara).Range.End)
    $delRange.Delete()
}

# --- 2. Drop the stray <w:lastRenderedPageBreak/> marker that precedes the
#        final "{registrationOn}" field run (a stale pagination artifact).
#        Re-typing the leading "{" character regenerates the run without it
#        while keeping the run's formatting (size/lang) intact.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$start = $lastPara.Range.Start
$r = $d.Range($start, $start + 1)
if ($r.Text -eq "{") {
    $r.Find.Execute("{", $true, $false, $false, $false, $false, $true, 1, $false, "{", 2)
}
